$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 788.06665
$ws.Range("J19").Value = 910.0833
$ws.Range("L19").Value = 910.0833
$ws.Range("N19").Value = -1260.0833

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 3099
$ws.Range("I51").Value = 2383.5
$ws.Range("J51").Value = 3351.5293
$ws.Range("K51").Value = 2383.5
$ws.Range("L51").Value = 3351.5293
$ws.Range("M51").Value = -1899.5
$ws.Range("N51").Value = -4319.5293

# Row 68 (Leve Item ID 10647)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71 (Leve Item ID 10647)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 97 (Leve Item ID 19885)
$ws.Range("H97").Value = 15249.25
$ws.Range("J97").Value = 15249.25
$ws.Range("L97").Value = 45747.75
$ws.Range("N97").Value = -46739.75

# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 1200
$ws.Range("J101").Value = 3871
$ws.Range("L101").Value = 11613
$ws.Range("N101").Value = -14857

# Row 133 (Leve Item ID 41856)
$ws.Range("H133").Value = 59999
$ws.Range("J133").Value = 59999
$ws.Range("L133").Value = 59999
$ws.Range("N133").Value = -70119

# Row 136 (Leve Item ID 42164)
$ws.Range("H136").Value = 211499.5
$ws.Range("J136").Value = 211499.5
$ws.Range("L136").Value = 211499.5
$ws.Range("N136").Value = -221699.5

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 6900.4443
$ws.Range("I137").Value = 7813.6665
$ws.Range("J137").Value = 2334.3333
$ws.Range("K137").Value = 23440.9995
$ws.Range("L137").Value = 7002.999899999999
$ws.Range("M137").Value = -20890.9995
$ws.Range("N137").Value = -12102.9999

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 6056.058
$ws.Range("I138").Value = 1259.9
$ws.Range("J138").Value = 6687.1313
$ws.Range("K138").Value = 3779.7
$ws.Range("L138").Value = 20061.3939
$ws.Range("M138").Value = 1360.3
$ws.Range("N138").Value = -30341.3939

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 2504.9075
$ws.Range("I32").Value = 2599.0833
$ws.Range("K32").Value = 2599.0833
$ws.Range("M32").Value = -2312.0833

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2843.0454
$ws.Range("I74").Value = 2835.5715
$ws.Range("K74").Value = 2835.5715
$ws.Range("M74").Value = -1961.5715

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2843.0454
$ws.Range("I77").Value = 2835.5715
$ws.Range("K77").Value = 14177.8575
$ws.Range("M77").Value = -9809.8575

# Row 133 (Leve Item ID 41857)
$ws.Range("H133").Value = 69664.836
$ws.Range("J133").Value = 69664.836
$ws.Range("L133").Value = 69664.836
$ws.Range("N133").Value = -74724.836

$ws = $wb.Worksheets.Item("BSM")
# Row 17 (Leve Item ID 2393)
$ws.Range("H17").Value = 5002.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5002.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5002.25
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5346.25

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1937.2222
$ws.Range("I94").Value = 1731.7576
$ws.Range("K94").Value = 1731.7576
$ws.Range("M94").Value = -1280.7576

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3979.3447
$ws.Range("I107").Value = 1685.8096
$ws.Range("K107").Value = 1685.8096
$ws.Range("M107").Value = 234.1904

# Row 132 (Leve Item ID 41855)
$ws.Range("H132").Value = 58852.668
$ws.Range("J132").Value = 58852.668
$ws.Range("L132").Value = 58852.668
$ws.Range("N132").Value = -68972.66800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 13 (Leve Item ID 1996)
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2075.7144
$ws.Range("I105").Value = 2503.8
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 2503.8
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = -756.8000000000002
$ws.Range("N105").Value = -4499.5

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (Leve Item ID 4867)
$ws.Range("H33").Value = 398.5
$ws.Range("I33").Value = 307
$ws.Range("K33").Value = 1842
$ws.Range("M33").Value = -1559

# Row 37 (Leve Item ID 9516)
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 41 (Leve Item ID 4700)
$ws.Range("H41").Value = 294
$ws.Range("I41").Value = 135.625
$ws.Range("J41").Value = 547.4
$ws.Range("K41").Value = 406.875
$ws.Range("L41").Value = 1642.2
$ws.Range("M41").Value = -68.875
$ws.Range("N41").Value = -2318.2

# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 2610.1875
$ws.Range("I132").Value = 2109.8572
$ws.Range("J132").Value = 2999.3333
$ws.Range("K132").Value = 18988.7148
$ws.Range("L132").Value = 26993.9997
$ws.Range("M132").Value = -16458.7148
$ws.Range("N132").Value = -32053.9997

# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 5248.76
$ws.Range("I138").Value = 2869.8333
$ws.Range("K138").Value = 8609.499899999999
$ws.Range("M138").Value = -3469.499899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 42 (Leve Item ID 27213)
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

# Row 115 (Leve Item ID 27213)
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()

# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 69999.25
$ws.Range("J123").Value = 69999.25
$ws.Range("L123").Value = 69999.25
$ws.Range("N123").Value = -74899.25

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3124.1667
$ws.Range("I126").Value = 3086.25
$ws.Range("K126").Value = 9258.75
$ws.Range("M126").Value = -6788.75

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2006.4
$ws.Range("J46").Value = 2316.3333
$ws.Range("L46").Value = 2316.3333
$ws.Range("N46").Value = -2692.3333

# Row 59 (Leve Item ID 25982)
$ws.Range("H59").Value = 25333
$ws.Range("J59").Value = 25333
$ws.Range("L59").Value = 25333
$ws.Range("N59").Value = -26641

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 1950.8
$ws.Range("I96").Value = 799.1429000000001
$ws.Range("K96").Value = 799.1429000000001
$ws.Range("M96").Value = 573.8570999999999

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 2498.6667
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 2248
$ws.Range("K113").Value = 9000
$ws.Range("L113").Value = 6744
$ws.Range("M113").Value = -6830
$ws.Range("N113").Value = -11084

# Row 114 (Leve Item ID 25978)
$ws.Range("H114").Value = 65130
$ws.Range("J114").Value = 65130
$ws.Range("L114").Value = 65130
$ws.Range("N114").Value = -73808

# Row 133 (Leve Item ID 41869)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
